$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22. This shifts the existing rows 22-42
# down to 23-43 (carrying their values/styles with them) and grows the
# sheet dimension to A1:R43 automatically.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with this week's new data point.
# (Same "Rabanito" / "Provincia de Chacabuco" / "Primera" pattern as the
# rest of the sheet, new date 2022-12-07 and new volume 7000.)
$ws.Cells.Item(22, 1).Value = 6
$ws.Cells.Item(22, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(22, 3).Value = "Metropolitana"
$ws.Cells.Item(22, 4).Value = 44902
$ws.Cells.Item(22, 5).Value = 13
$ws.Cells.Item(22, 6).Value = 300000001
$ws.Cells.Item(22, 7).Value = "Rabanito"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 7000
$ws.Cells.Item(22, 11).Value = 3000
$ws.Cells.Item(22, 12).Value = 3000
$ws.Cells.Item(22, 13).Value = 3000
$ws.Cells.Item(22, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(22, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(22, 16).Value = 30
$ws.Cells.Item(22, 17).Value = 100
$ws.Cells.Item(22, 18).Value = "Hortaliza"
